# Appends the new 2023-2024 Croatia HNL fixtures (rows 82-93) scraped on 19-12-2023
# to the existing results table on Sheet1, matching the formatting of the last
# pre-existing row (bold/bordered "Indice" column, datetime-formatted match column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 81
$firstNewRow = 82
$lastNewRow = 93

# Seed the new rows formatting (column styles only) by copying the A:V slice of the
# last existing row down across the whole new block in one shot.
$srcRange = $ws.Range("A" + $lastExistingRow + ":V" + $lastExistingRow)
$dstRange = $ws.Range("A" + $firstNewRow + ":V" + $lastNewRow)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 82 (Indice 81)
$ws.Range("A82").Value2 = 81
$ws.Range("B82").Value2 = "croatia"
$ws.Range("C82").Value2 = "hnl"
$ws.Range("D82").Value2 = "2023-2024"
$ws.Range("E82").Value2 = 45263.625
$ws.Range("F82").Value2 = "Rudes"
$ws.Range("G82").Value2 = 0
$ws.Range("H82").Value2 = "Lok. Zagreb"
$ws.Range("I82").Value2 = 0
$ws.Range("J82").Value2 = 2.96
$ws.Range("K82").Value2 = "26/11/2023 17:12"
$ws.Range("L82").Value2 = 5.16
$ws.Range("M82").Value2 = "03/12/2023 14:57"
$ws.Range("N82").Value2 = 3.3
$ws.Range("O82").Value2 = "26/11/2023 17:12"
$ws.Range("P82").Value2 = 3.96
$ws.Range("Q82").Value2 = "03/12/2023 14:57"
$ws.Range("R82").Value2 = 2.42
$ws.Range("S82").Value2 = "26/11/2023 17:12"
$ws.Range("T82").Value2 = 1.66
$ws.Range("U82").Value2 = "03/12/2023 14:57"
$ws.Range("V82").Value2 = "https://www.betexplorer.com/football/croatia/hnl/rudes-lok-zagreb/vqqlP8Xj/"

# Row 83 (Indice 82)
$ws.Range("A83").Value2 = 82
$ws.Range("B83").Value2 = "croatia"
$ws.Range("C83").Value2 = "hnl"
$ws.Range("D83").Value2 = "2023-2024"
$ws.Range("E83").Value2 = 45264.70833333334
$ws.Range("F83").Value2 = "Slaven Belupo"
$ws.Range("G83").Value2 = 0
$ws.Range("H83").Value2 = "D. Zagreb"
$ws.Range("I83").Value2 = 2
$ws.Range("J83").Value2 = 7.38
$ws.Range("K83").Value2 = "28/11/2023 12:12"
$ws.Range("L83").Value2 = 7.89
$ws.Range("M83").Value2 = "04/12/2023 16:59"
$ws.Range("N83").Value2 = 4.28
$ws.Range("O83").Value2 = "28/11/2023 12:12"
$ws.Range("P83").Value2 = 4.62
$ws.Range("Q83").Value2 = "04/12/2023 16:59"
$ws.Range("R83").Value2 = 1.45
$ws.Range("S83").Value2 = "28/11/2023 12:12"
$ws.Range("T83").Value2 = 1.41
$ws.Range("U83").Value2 = "04/12/2023 16:59"
$ws.Range("V83").Value2 = "https://www.betexplorer.com/football/croatia/hnl/slaven-belupo-din-zagreb/6kXVmhAi/"

# Row 84 (Indice 83)
$ws.Range("A84").Value2 = 83
$ws.Range("B84").Value2 = "croatia"
$ws.Range("C84").Value2 = "hnl"
$ws.Range("D84").Value2 = "2023-2024"
$ws.Range("E84").Value2 = 45268.72916666666
$ws.Range("F84").Value2 = "Lok. Zagreb"
$ws.Range("G84").Value2 = 1
$ws.Range("H84").Value2 = "Hajduk Split"
$ws.Range("I84").Value2 = 1
$ws.Range("J84").Value2 = 4.21
$ws.Range("K84").Value2 = "03/12/2023 15:12"
$ws.Range("L84").Value2 = 4.61
$ws.Range("M84").Value2 = "08/12/2023 17:29"
$ws.Range("N84").Value2 = 3.63
$ws.Range("O84").Value2 = "03/12/2023 15:12"
$ws.Range("P84").Value2 = 3.66
$ws.Range("Q84").Value2 = "08/12/2023 17:29"
$ws.Range("R84").Value2 = 1.78
$ws.Range("S84").Value2 = "03/12/2023 15:12"
$ws.Range("T84").Value2 = 1.79
$ws.Range("U84").Value2 = "08/12/2023 17:29"
$ws.Range("V84").Value2 = "https://www.betexplorer.com/football/croatia/hnl/lok-zagreb-hajduk-split/ELNjqUOM/"

# Row 85 (Indice 84)
$ws.Range("A85").Value2 = 84
$ws.Range("B85").Value2 = "croatia"
$ws.Range("C85").Value2 = "hnl"
$ws.Range("D85").Value2 = "2023-2024"
$ws.Range("E85").Value2 = 45269.61805555555
$ws.Range("F85").Value2 = "Gorica"
$ws.Range("G85").Value2 = 0
$ws.Range("H85").Value2 = "Istra 1961"
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 1.74
$ws.Range("K85").Value2 = "02/12/2023 21:42"
$ws.Range("L85").Value2 = 1.99
$ws.Range("M85").Value2 = "09/12/2023 14:49"
$ws.Range("N85").Value2 = 3.53
$ws.Range("O85").Value2 = "02/12/2023 21:42"
$ws.Range("P85").Value2 = 3.26
$ws.Range("Q85").Value2 = "09/12/2023 14:45"
$ws.Range("R85").Value2 = 4.53
$ws.Range("S85").Value2 = "02/12/2023 21:42"
$ws.Range("T85").Value2 = 4.24
$ws.Range("U85").Value2 = "09/12/2023 14:49"
$ws.Range("V85").Value2 = "https://www.betexplorer.com/football/croatia/hnl/hnk-gorica-istra-1961/WdIerlvT/"

# Row 86 (Indice 85)
$ws.Range("A86").Value2 = 85
$ws.Range("B86").Value2 = "croatia"
$ws.Range("C86").Value2 = "hnl"
$ws.Range("D86").Value2 = "2023-2024"
$ws.Range("E86").Value2 = 45269.70833333334
$ws.Range("F86").Value2 = "D. Zagreb"
$ws.Range("G86").Value2 = 1
$ws.Range("H86").Value2 = "Rudes"
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 1.13
$ws.Range("K86").Value2 = "04/12/2023 17:12"
$ws.Range("L86").Value2 = 1.1
$ws.Range("M86").Value2 = "09/12/2023 16:56"
$ws.Range("N86").Value2 = 8.19
$ws.Range("O86").Value2 = "04/12/2023 17:12"
$ws.Range("P86").Value2 = 9.970000000000001
$ws.Range("Q86").Value2 = "09/12/2023 16:56"
$ws.Range("R86").Value2 = 13.9
$ws.Range("S86").Value2 = "04/12/2023 17:12"
$ws.Range("T86").Value2 = 24.38
$ws.Range("U86").Value2 = "09/12/2023 16:56"
$ws.Range("V86").Value2 = "https://www.betexplorer.com/football/croatia/hnl/din-zagreb-rudes/zVOnpA9G/"

# Row 87 (Indice 86)
$ws.Range("A87").Value2 = 86
$ws.Range("B87").Value2 = "croatia"
$ws.Range("C87").Value2 = "hnl"
$ws.Range("D87").Value2 = "2023-2024"
$ws.Range("E87").Value2 = 45270.625
$ws.Range("F87").Value2 = "Varazdin"
$ws.Range("G87").Value2 = 2
$ws.Range("H87").Value2 = "Osijek"
$ws.Range("I87").Value2 = 2
$ws.Range("J87").Value2 = 3.18
$ws.Range("K87").Value2 = "03/12/2023 15:12"
$ws.Range("L87").Value2 = 3.37
$ws.Range("M87").Value2 = "10/12/2023 14:56"
$ws.Range("N87").Value2 = 3.34
$ws.Range("O87").Value2 = "03/12/2023 15:12"
$ws.Range("P87").Value2 = 3.56
$ws.Range("Q87").Value2 = "10/12/2023 14:56"
$ws.Range("R87").Value2 = 2.17
$ws.Range("S87").Value2 = "03/12/2023 15:12"
$ws.Range("T87").Value2 = 2.14
$ws.Range("U87").Value2 = "10/12/2023 14:56"
$ws.Range("V87").Value2 = "https://www.betexplorer.com/football/croatia/hnl/varazdin-osijek/S4LvnWv4/"

# Row 88 (Indice 87)
$ws.Range("A88").Value2 = 87
$ws.Range("B88").Value2 = "croatia"
$ws.Range("C88").Value2 = "hnl"
$ws.Range("D88").Value2 = "2023-2024"
$ws.Range("E88").Value2 = 45270.71527777778
$ws.Range("F88").Value2 = "Rijeka"
$ws.Range("G88").Value2 = 2
$ws.Range("H88").Value2 = "Slaven Belupo"
$ws.Range("I88").Value2 = 4
$ws.Range("J88").Value2 = 1.25
$ws.Range("K88").Value2 = "04/12/2023 17:12"
$ws.Range("L88").Value2 = 1.26
$ws.Range("M88").Value2 = "10/12/2023 17:05"
$ws.Range("N88").Value2 = 5.55
$ws.Range("O88").Value2 = "04/12/2023 17:12"
$ws.Range("P88").Value2 = 5.77
$ws.Range("Q88").Value2 = "10/12/2023 17:05"
$ws.Range("R88").Value2 = 9.109999999999999
$ws.Range("S88").Value2 = "04/12/2023 17:12"
$ws.Range("T88").Value2 = 11.57
$ws.Range("U88").Value2 = "10/12/2023 17:05"
$ws.Range("V88").Value2 = "https://www.betexplorer.com/football/croatia/hnl/rijeka-slaven-belupo/t8ProjfA/"

# Row 89 (Indice 88)
$ws.Range("A89").Value2 = 88
$ws.Range("B89").Value2 = "croatia"
$ws.Range("C89").Value2 = "hnl"
$ws.Range("D89").Value2 = "2023-2024"
$ws.Range("E89").Value2 = 45275.70833333334
$ws.Range("F89").Value2 = "Gorica"
$ws.Range("G89").Value2 = 1
$ws.Range("H89").Value2 = "Varazdin"
$ws.Range("I89").Value2 = 3
$ws.Range("J89").Value2 = 1.95
$ws.Range("K89").Value2 = "10/12/2023 15:13"
$ws.Range("L89").Value2 = 2.41
$ws.Range("M89").Value2 = "15/12/2023 16:58"
$ws.Range("N89").Value2 = 3.32
$ws.Range("O89").Value2 = "10/12/2023 15:13"
$ws.Range("P89").Value2 = 3.23
$ws.Range("Q89").Value2 = "15/12/2023 16:46"
$ws.Range("R89").Value2 = 3.83
$ws.Range("S89").Value2 = "10/12/2023 15:13"
$ws.Range("T89").Value2 = 3.13
$ws.Range("U89").Value2 = "15/12/2023 16:58"
$ws.Range("V89").Value2 = "https://www.betexplorer.com/football/croatia/hnl/hnk-gorica-varazdin/Y54Nw8np/"

# Row 90 (Indice 89)
$ws.Range("A90").Value2 = 89
$ws.Range("B90").Value2 = "croatia"
$ws.Range("C90").Value2 = "hnl"
$ws.Range("D90").Value2 = "2023-2024"
$ws.Range("E90").Value2 = 45276.61805555555
$ws.Range("F90").Value2 = "Lok. Zagreb"
$ws.Range("G90").Value2 = 3
$ws.Range("H90").Value2 = "Istra 1961"
$ws.Range("I90").Value2 = 0
$ws.Range("J90").Value2 = 1.76
$ws.Range("K90").Value2 = "09/12/2023 15:13"
$ws.Range("L90").Value2 = 1.98
$ws.Range("M90").Value2 = "16/12/2023 14:45"
$ws.Range("N90").Value2 = 3.56
$ws.Range("O90").Value2 = "09/12/2023 15:13"
$ws.Range("P90").Value2 = 3.47
$ws.Range("Q90").Value2 = "16/12/2023 14:46"
$ws.Range("R90").Value2 = 4.36
$ws.Range("S90").Value2 = "09/12/2023 15:13"
$ws.Range("T90").Value2 = 3.94
$ws.Range("U90").Value2 = "16/12/2023 14:45"
$ws.Range("V90").Value2 = "https://www.betexplorer.com/football/croatia/hnl/lok-zagreb-istra-1961/SS3RxS1j/"

# Row 91 (Indice 90)
$ws.Range("A91").Value2 = 90
$ws.Range("B91").Value2 = "croatia"
$ws.Range("C91").Value2 = "hnl"
$ws.Range("D91").Value2 = "2023-2024"
$ws.Range("E91").Value2 = 45276.70833333334
$ws.Range("F91").Value2 = "Rijeka"
$ws.Range("G91").Value2 = 3
$ws.Range("H91").Value2 = "Rudes"
$ws.Range("I91").Value2 = 0
$ws.Range("J91").Value2 = 1.18
$ws.Range("K91").Value2 = "10/12/2023 17:12"
$ws.Range("L91").Value2 = 1.18
$ws.Range("M91").Value2 = "16/12/2023 16:49"
$ws.Range("N91").Value2 = 6.47
$ws.Range("O91").Value2 = "10/12/2023 17:12"
$ws.Range("P91").Value2 = 7.19
$ws.Range("Q91").Value2 = "16/12/2023 16:58"
$ws.Range("R91").Value2 = 10.75
$ws.Range("S91").Value2 = "10/12/2023 17:12"
$ws.Range("T91").Value2 = 14.14
$ws.Range("U91").Value2 = "16/12/2023 16:58"
$ws.Range("V91").Value2 = "https://www.betexplorer.com/football/croatia/hnl/rijeka-rudes/ryYBpvLd/"

# Row 92 (Indice 91)
$ws.Range("A92").Value2 = 91
$ws.Range("B92").Value2 = "croatia"
$ws.Range("C92").Value2 = "hnl"
$ws.Range("D92").Value2 = "2023-2024"
$ws.Range("E92").Value2 = 45277.60416666666
$ws.Range("F92").Value2 = "Osijek"
$ws.Range("G92").Value2 = 4
$ws.Range("H92").Value2 = "Slaven Belupo"
$ws.Range("I92").Value2 = 1
$ws.Range("J92").Value2 = 1.5
$ws.Range("K92").Value2 = "10/12/2023 17:12"
$ws.Range("L92").Value2 = 1.57
$ws.Range("M92").Value2 = "17/12/2023 14:26"
$ws.Range("N92").Value2 = 4.14
$ws.Range("O92").Value2 = "10/12/2023 17:12"
$ws.Range("P92").Value2 = 3.9
$ws.Range("Q92").Value2 = "17/12/2023 14:26"
$ws.Range("R92").Value2 = 5.7
$ws.Range("S92").Value2 = "10/12/2023 17:12"
$ws.Range("T92").Value2 = 6.41
$ws.Range("U92").Value2 = "17/12/2023 14:26"
$ws.Range("V92").Value2 = "https://www.betexplorer.com/football/croatia/hnl/osijek-slaven-belupo/MRWFqbz3/"

# Row 93 (Indice 92)
$ws.Range("A93").Value2 = 92
$ws.Range("B93").Value2 = "croatia"
$ws.Range("C93").Value2 = "hnl"
$ws.Range("D93").Value2 = "2023-2024"
$ws.Range("E93").Value2 = 45277.71875
$ws.Range("F93").Value2 = "D. Zagreb"
$ws.Range("G93").Value2 = 0
$ws.Range("H93").Value2 = "Hajduk Split"
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 1.86
$ws.Range("K93").Value2 = "10/12/2023 19:12"
$ws.Range("L93").Value2 = 1.91
$ws.Range("M93").Value2 = "17/12/2023 17:14"
$ws.Range("N93").Value2 = 3.44
$ws.Range("O93").Value2 = "10/12/2023 19:12"
$ws.Range("P93").Value2 = 3.53
$ws.Range("Q93").Value2 = "17/12/2023 17:14"
$ws.Range("R93").Value2 = 3.99
$ws.Range("S93").Value2 = "10/12/2023 19:12"
$ws.Range("T93").Value2 = 4.16
$ws.Range("U93").Value2 = "17/12/2023 17:14"
$ws.Range("V93").Value2 = "https://www.betexplorer.com/football/croatia/hnl/din-zagreb-hajduk-split/v37VynHd/"

